# Update "paises.xlsx" COVID dashboard with refreshed country data
# (commit: "Update countries & provincias Spain")
#
# The sheet is sorted by "Casos totales" (col B) descending. Refreshing the
# source numbers nudges a handful of countries past their neighbours, so a
# few rows swap identities in addition to getting new B:H figures:
#   - Kuwait (row 55) overtakes Finlandia (row 54)               -> rows 54/55 swap
#   - Yemen's new total (21) jumps it above Granada and six more
#     small countries that were tied/close, inserting it at row 186
#     and pushing Granada..Islas Malvinas down one row each
#     (the previous last row of that block, old "Yemen" row, is
#     absorbed by the shift).
# Everything else in the sheet is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Mayo de 2020 a las 14:03"

# --- Straightforward numeric refreshes (country stays put) ------------
# Row 25: Suecia
$ws.Cells.Item(25, 2).Value = 23216
$ws.Cells.Item(25, 3).Value = 495
$ws.Cells.Item(25, 5).Value = 16288
$ws.Cells.Item(25, 6).Value = 435
$ws.Cells.Item(25, 7).Value = 85
$ws.Cells.Item(25, 8).Value = 2854

# Row 60: Kazajistan
$ws.Cells.Item(60, 2).Value = 4179
$ws.Cells.Item(60, 3).Value = 130
$ws.Cells.Item(60, 5).Value = 2886

# Row 89: Senegal
$ws.Cells.Item(89, 2).Value = 1329
$ws.Cells.Item(89, 3).Value = 58
$ws.Cells.Item(89, 4).Value = 470
$ws.Cells.Item(89, 5).Value = 848
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 11

# --- Kuwait / Finlandia swap (rows 54-55) ------------------------------
# Row 54 becomes Kuwait with its updated figures
$ws.Cells.Item(54, 1).Value = "Kuwait"
$ws.Cells.Item(54, 2).Value = 5804
$ws.Cells.Item(54, 3).Value = 526
$ws.Cells.Item(54, 4).Value = 2032
$ws.Cells.Item(54, 5).Value = 3732
$ws.Cells.Item(54, 6).Value = 90
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 40

# Row 55 becomes Finlandia (its old figures, unchanged from row 54 before)
$ws.Cells.Item(55, 1).Value = "Finlandia"
$ws.Cells.Item(55, 2).Value = 5412
$ws.Cells.Item(55, 3).Value = 85
$ws.Cells.Item(55, 4).Value = 3500
$ws.Cells.Item(55, 5).Value = 1672
$ws.Cells.Item(55, 6).Value = 49
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 240

# --- Yemen inserted ahead of Granada; rows 186-202 shift down one ------
# Row 186: new Yemen entry (brand-new figures)
$ws.Cells.Item(186, 1).Value = "Yemen"
$ws.Cells.Item(186, 2).Value = 21
$ws.Cells.Item(186, 3).Value = 9
$ws.Cells.Item(186, 4).Value = 1
$ws.Cells.Item(186, 5).Value = 17
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 1
$ws.Cells.Item(186, 8).Value = 3

# Row 187: Granada (was row 186's data)
$ws.Cells.Item(187, 1).Value = "Granada"
$ws.Cells.Item(187, 2).Value = 21
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 13
$ws.Cells.Item(187, 5).Value = 8
$ws.Cells.Item(187, 6).Value = 4
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

# Row 188: Laos (was row 187's data)
$ws.Cells.Item(188, 1).Value = "Laos"
$ws.Cells.Item(188, 2).Value = 19
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 9
$ws.Cells.Item(188, 5).Value = 10
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

# Row 189: Fiyi (was row 188's data)
$ws.Cells.Item(189, 1).Value = "Fiyi"
$ws.Cells.Item(189, 2).Value = 18
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 14
$ws.Cells.Item(189, 5).Value = 4
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0

# Row 190: Santa Lucia (was row 189's data)
$ws.Cells.Item(190, 1).Value = "Santa Lucia"
$ws.Cells.Item(190, 2).Value = 18
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 15
$ws.Cells.Item(190, 5).Value = 3
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# Row 191: Belice (was row 190's data)
$ws.Cells.Item(191, 1).Value = "Belice"
$ws.Cells.Item(191, 2).Value = 18
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 14
$ws.Cells.Item(191, 5).Value = 2
$ws.Cells.Item(191, 6).Value = 1
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 2

# Row 192: Nueva Caledonia (was row 191's data)
$ws.Cells.Item(192, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(192, 2).Value = 18
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 17
$ws.Cells.Item(192, 5).Value = 1
$ws.Cells.Item(192, 6).Value = 1
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

# Row 193: Islas Virgenes de los Estados Unidos (was row 192's data)
$ws.Cells.Item(193, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(193, 2).Value = 17
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 5).Value = 17
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

# Row 194: San Vicente y las Granadinas (was row 193's data)
$ws.Cells.Item(194, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(194, 2).Value = 17
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 9
$ws.Cells.Item(194, 5).Value = 8
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

# Row 195: Gambia (was row 194's data)
$ws.Cells.Item(195, 1).Value = "Gambia"
$ws.Cells.Item(195, 2).Value = 17
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 9
$ws.Cells.Item(195, 5).Value = 7
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 1

# Row 196: Namibia (was row 195's data)
$ws.Cells.Item(196, 1).Value = "Namibia"
$ws.Cells.Item(196, 2).Value = 16
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 8
$ws.Cells.Item(196, 5).Value = 8
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0

# Row 197: Dominica (was row 196's data)
$ws.Cells.Item(197, 1).Value = "Dominica"
$ws.Cells.Item(197, 2).Value = 16
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 13
$ws.Cells.Item(197, 5).Value = 3
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0

# Row 198: Curazao (was row 197's data)
$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 2).Value = 16
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 13
$ws.Cells.Item(198, 5).Value = 2
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 1

# Row 199: San Cristobal y Nieves (was row 198's data)
$ws.Cells.Item(199, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(199, 2).Value = 15
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 8
$ws.Cells.Item(199, 5).Value = 7
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# Row 200: Burundi (was row 199's data)
$ws.Cells.Item(200, 1).Value = "Burundi"
$ws.Cells.Item(200, 2).Value = 15
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 7
$ws.Cells.Item(200, 5).Value = 7
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 1

# Row 201: Nicaragua (was row 200's data)
$ws.Cells.Item(201, 1).Value = "Nicaragua"
$ws.Cells.Item(201, 2).Value = 15
$ws.Cells.Item(201, 3).Value = 0
$ws.Cells.Item(201, 4).Value = 7
$ws.Cells.Item(201, 5).Value = 3
$ws.Cells.Item(201, 6).Value = 0
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 5

# Row 202: Islas Malvinas (was row 201's data) - absorbs the old Yemen slot
$ws.Cells.Item(202, 1).Value = "Islas Malvinas"
$ws.Cells.Item(202, 2).Value = 13
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 13
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203 (Islas Turcas y Caicos) and everything below is unchanged.
